$d = $word.ActiveDocument

# --- 1. "La agregación representa a los partidos " -> insert "más interna "
#        after "La agregación " (paragraph stays a single list item; text only
#        changes).
$d.Content.Find.Execute(
    "La agregación representa a los partidos",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "La agregación más interna representa a los partidos",
    2) | Out-Null

# --- 2. "Decisiones en cuanto al usuario" becomes a new bulleted/numbered
#        list item (same list as the paragraph above, numId=3) with new text.
$pModel = $d.Paragraphs.Item(2)
$pDecisiones = $d.Paragraphs.Item(3)
$pDecisiones.Range.Text = "La segunda agregación externa es con propósito de saber de que campeonato el alumno esta prediciendo el equipo campeón y subcampeón."
$pDecisiones.Style = "Prrafodelista"
$pDecisiones.Range.ListFormat.ListId = $pModel.Range.ListFormat.ListId
$pDecisiones.Range.ListFormat.ListLevelNumber = $pModel.Range.ListFormat.ListLevelNumber

# --- 3. "Se utilizará "ci" como identificador..." paragraph: move it onto the
#        numId=3 list, merge the three runs (dropping the spell-check
#        proofErr bookmarks) into a single run with the same visible text.
$pCi = $d.Paragraphs.Item(4)
$rngCi = $pCi.Range
$rngCi.MoveEnd(1, -1) | Out-Null
$rngCi.Delete()
$rngCi.InsertAfter("Se utilizará “ci” como identificador ya que de esta forma evitaremos que una misma persona se cree más de una cuenta en el sistema.")
$pCi.Range.ListFormat.ListId = $pModel.Range.ListFormat.ListId
$pCi.Range.ListFormat.ListLevelNumber = $pModel.Range.ListFormat.ListLevelNumber

# --- 4. Two brand new list items right after the "ci" paragraph, on the same
#        numId=3 list.
$pCi.Range.InsertParagraphAfter()
$pCarrera = $d.Paragraphs.Item(5)
$pCarrera.Range.InsertBefore("Un alumno pertenece a por lo menos a una carrera.")

$pCarrera.Range.InsertParagraphAfter()
$pEquipo = $d.Paragraphs.Item(6)
$pEquipo.Range.InsertBefore("Un equipo por lo menos juega un partido ")
